$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# DMSEmailTitle value (B23): drop the space before the first placeholder
$ws.Range("B23").Value = "Compliance Waste Returns{0} {1} - Email and Submission"

# DMSExcelReturnTitle value (B25): drop the space and renumber the placeholders
$ws.Range("B25").Value = "Compliance Waste Returns{0} {1}"

# Update the active selection to B22 (was B25)
$ws.Range("B22").Select()
